# Uploaded from Matlab using gitup()
# Adds eeg_lbl_path hyperlinks (col K) and have_seizure_lbl flags (col J)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window position (matches the saved workbook view) ---
$wb.Windows.Item(1).Left = 8940

# --- eeg_lbl_path column (K), header + hyperlinks to UNC paths ---
$ws.Range("K1").Value = "eeg_lbl_path"
$ws.Hyperlinks.Add($ws.Range("K2"), "\\neurodata\Lab Neurophysiology root\EEG conversion\NatySST_TdTET339")
$ws.Hyperlinks.Add($ws.Range("K4"), "\\neurodata\Lab Neurophysiology root\EEG conversion\NatymTORET283")
$ws.Hyperlinks.Add($ws.Range("K5"), "\\neurodata\Lab Neurophysiology root\EEG Naty\mTOR MUT\Naty SST_TdT ET 343")
$ws.Hyperlinks.Add($ws.Range("K6"), "\\neurodata\Lab Neurophysiology root\EEG Naty\mTOR MUT\Naty SST_TdT ET 413")

# --- have_seizure_lbl column (J) ---
$ws.Range("J1").Value = "have_seizure_lbl"
$ws.Range("J2").Value = 1
$ws.Range("J3").Value = 0
$ws.Range("J4").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("J7").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("J13").Value = 0

# --- Column widths for the new columns ---
$ws.Columns.Item(10).ColumnWidth = 15.5833333333
$ws.Columns.Item(11).ColumnWidth = 31.4166666667

# --- Selection state, matching the saved view ---
$ws.Range("K2").Select()

Write-Host "done"
